$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "MuSCs"
$ws.Range("G2").Value = 5.791393
$ws.Range("H2").Value = 17.374179
$ws.Range("I2").Value = 0.2508946350719245
$ws.Range("J2").Value = 0.2508946350719244
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.690195
$ws.Range("N2").Value = 2.070585
$ws.Range("O2").Value = 0.1060250152438306
$ws.Range("P2").Value = 0.1060250152438306
$ws.Range("Q2").Value = 3.997190491635
$ws.Range("R2").Value = 35.974714424715
$ws.Range("S2").Value = 0.02660110750809612
$ws.Range("T2").Value = 0.02660110750809611
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("G3").Value = 5.791393
$ws.Range("H3").Value = 17.374179
$ws.Range("I3").Value = 0.2508946350719245
$ws.Range("J3").Value = 0.2508946350719244
$ws.Range("M3").Value = 5.819542333333334
$ws.Range("N3").Value = 17.458627
$ws.Range("O3").Value = 0.8939749847561693
$ws.Range("P3").Value = 0.8939749847561693
$ws.Range("Q3").Value = 33.70325673247034
$ws.Range("R3").Value = 303.329310592233
$ws.Range("S3").Value = 0.2242935275638284
$ws.Range("T3").Value = 0.2242935275638283
$ws.Range("A4").Value = "FAPs"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("G4").Value = 8.502533
$ws.Range("H4").Value = 25.507599
$ws.Range("I4").Value = 0.3683465988617928
$ws.Range("J4").Value = 0.3683465988617928
$ws.Range("M4").Value = 0.690195
$ws.Range("N4").Value = 2.070585
$ws.Range("O4").Value = 0.1060250152438306
$ws.Range("P4").Value = 0.1060250152438306
$ws.Range("Q4").Value = 5.868405763935
$ws.Range("R4").Value = 52.81565187541499
$ws.Range("S4").Value = 0.03905395375933475
$ws.Range("T4").Value = 0.03905395375933474
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("I5").Value = 0.3683465988617928
$ws.Range("J5").Value = 0.3683465988617928
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 5.819542333333334
$ws.Range("N5").Value = 17.458627
$ws.Range("O5").Value = 0.8939749847561693
$ws.Range("P5").Value = 0.8939749847561693
$ws.Range("Q5").Value = 49.48085073406367
$ws.Range("R5").Value = 445.327656606573
$ws.Range("S5").Value = 0.329292645102458
$ws.Range("T5").Value = 0.329292645102458
$ws.Range("A6").Value = "MuSCs"
$ws.Range("G6").Value = 2.834746
$ws.Range("H6").Value = 8.504238
$ws.Range("I6").Value = 0.1228068209481894
$ws.Range("J6").Value = 0.1228068209481894
$ws.Range("M6").Value = 0.690195
$ws.Range("N6").Value = 2.070585
$ws.Range("O6").Value = 0.1060250152438306
$ws.Range("P6").Value = 0.1060250152438306
$ws.Range("Q6").Value = 1.95652751547
$ws.Range("R6").Value = 17.60874763923
$ws.Range("S6").Value = 0.01302059506307816
$ws.Range("T6").Value = 0.01302059506307816
$ws.Range("A7").Value = "MuSCs"
$ws.Range("G7").Value = 2.834746
$ws.Range("H7").Value = 8.504238
$ws.Range("I7").Value = 0.1228068209481894
$ws.Range("J7").Value = 0.1228068209481894
$ws.Range("M7").Value = 5.819542333333334
$ws.Range("N7").Value = 17.458627
$ws.Range("O7").Value = 0.8939749847561693
$ws.Range("P7").Value = 0.8939749847561693
$ws.Range("Q7").Value = 16.49692435124734
$ws.Range("R7").Value = 148.472319161226
$ws.Range("S7").Value = 0.1097862258851113
$ws.Range("T7").Value = 0.1097862258851113
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("G8").Value = 5.954296666666667
$ws.Range("H8").Value = 17.86289
$ws.Range("I8").Value = 0.2579519451180933
$ws.Range("J8").Value = 0.2579519451180933
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.690195
$ws.Range("N8").Value = 2.070585
$ws.Range("O8").Value = 0.1060250152438306
$ws.Range("P8").Value = 0.1060250152438306
$ws.Range("Q8").Value = 4.109625787850001
$ws.Range("R8").Value = 36.98663209065
$ws.Range("S8").Value = 0.0273493589133216
$ws.Range("T8").Value = 0.0273493589133216
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("G9").Value = 5.954296666666667
$ws.Range("H9").Value = 17.86289
$ws.Range("I9").Value = 0.2579519451180933
$ws.Range("J9").Value = 0.2579519451180933
$ws.Range("M9").Value = 5.819542333333334
$ws.Range("N9").Value = 17.458627
$ws.Range("O9").Value = 0.8939749847561693
$ws.Range("P9").Value = 0.8939749847561693
$ws.Range("Q9").Value = 34.65128151689223
$ws.Range("R9").Value = 311.86153365203
$ws.Range("S9").Value = 0.2306025862047717
$ws.Range("T9").Value = 0.2306025862047717
# Remove the now-superseded rows 10-13 (data fully covered by rows 2-9 after TPM update)
$ws.Rows("10:13").Delete()
